$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F-column "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1070
$ws1.Range("F5").Value = 2921
$ws1.Range("F6").Value = 97
$ws1.Range("F7").Value = 289
$ws1.Range("F8").Value = 31
$ws1.Range("F11").Value = 104
$ws1.Range("F12").Value = 153
$ws1.Range("F13").Value = 69
$ws1.Range("F14").Value = 2761
$ws1.Range("F15").Value = 1030

# Sheet "全部类型" (sheet4): F-column "想去人数" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1070
$ws4.Range("F6").Value = 2921
$ws4.Range("F7").Value = 97
$ws4.Range("F8").Value = 289
$ws4.Range("F9").Value = 31
$ws4.Range("F13").Value = 104
$ws4.Range("F14").Value = 153
$ws4.Range("F15").Value = 69
$ws4.Range("F16").Value = 2761
$ws4.Range("F17").Value = 1030
